$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 155.2138263333333
$ws.Range("H2").Value = 465.641479
$ws.Range("I2").Value = 0.3492508712612995
$ws.Range("J2").Value = 0.3492508712612995
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 153.5290173333333
$ws.Range("N2").Value = 460.587052
$ws.Range("O2").Value = 0.3172206968818489
$ws.Range("P2").Value = 0.317220696881849
$ws.Range("Q2").Value = 23829.82623350332
$ws.Range("R2").Value = 214468.4361015299
$ws.Range("S2").Value = 0.1107896047681023
$ws.Range("T2").Value = 0.1107896047681023

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 155.2138263333333
$ws.Range("H3").Value = 465.641479
$ws.Range("I3").Value = 0.3492508712612995
$ws.Range("J3").Value = 0.3492508712612995
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 168.7997026666667
$ws.Range("N3").Value = 506.3991080000001
$ws.Range("O3").Value = 0.3487728915577651
$ws.Range("P3").Value = 0.3487728915577651
$ws.Range("Q3").Value = 26200.04773482231
$ws.Range("R3").Value = 235800.4296134008
$ws.Range("S3").Value = 0.1218092362488722
$ws.Range("T3").Value = 0.1218092362488722

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 155.2138263333333
$ws.Range("H4").Value = 465.641479
$ws.Range("I4").Value = 0.3492508712612995
$ws.Range("J4").Value = 0.3492508712612995
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 68.09032333333333
$ws.Range("N4").Value = 204.27097
$ws.Range("O4").Value = 0.1406878008722904
$ws.Range("P4").Value = 0.1406878008722904
$ws.Range("Q4").Value = 10568.55962084051
$ws.Range("R4").Value = 95117.03658756464
$ws.Range("S4").Value = 0.04913533703048363
$ws.Range("T4").Value = 0.04913533703048364

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 155.2138263333333
$ws.Range("H5").Value = 465.641479
$ws.Range("I5").Value = 0.3492508712612995
$ws.Range("J5").Value = 0.3492508712612995
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 93.562673
$ws.Range("N5").Value = 280.688019
$ws.Range("O5").Value = 0.1933186106880956
$ws.Range("P5").Value = 0.1933186106880956
$ws.Range("Q5").Value = 14522.22047830446
$ws.Range("R5").Value = 130699.9843047401
$ws.Range("S5").Value = 0.06751669321384136
$ws.Range("T5").Value = 0.06751669321384136

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 246.1811623333333
$ws.Range("H6").Value = 738.543487
$ws.Range("I6").Value = 0.5539389593320749
$ws.Range("J6").Value = 0.5539389593320749
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 153.5290173333333
$ws.Range("N6").Value = 460.587052
$ws.Range("O6").Value = 0.3172206968818489
$ws.Range("P6").Value = 0.317220696881849
$ws.Range("Q6").Value = 37795.95193901448
$ws.Range("R6").Value = 340163.5674511303
$ws.Range("S6").Value = 0.175720902709327
$ws.Range("T6").Value = 0.175720902709327

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 246.1811623333333
$ws.Range("H7").Value = 738.543487
$ws.Range("I7").Value = 0.5539389593320749
$ws.Range("J7").Value = 0.5539389593320749
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 168.7997026666667
$ws.Range("N7").Value = 506.3991080000001
$ws.Range("O7").Value = 0.3487728915577651
$ws.Range("P7").Value = 0.3487728915577651
$ws.Range("Q7").Value = 41555.30700400107
$ws.Range("R7").Value = 373997.7630360097
$ws.Range("S7").Value = 0.193198892592747
$ws.Range("T7").Value = 0.193198892592747

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 246.1811623333333
$ws.Range("H8").Value = 738.543487
$ws.Range("I8").Value = 0.5539389593320749
$ws.Range("J8").Value = 0.5539389593320749
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 68.09032333333333
$ws.Range("N8").Value = 204.27097
$ws.Range("O8").Value = 0.1406878008722904
$ws.Range("P8").Value = 0.1406878008722904
$ws.Range("Q8").Value = 16762.55494185249
$ws.Range("R8").Value = 150862.9944766724
$ws.Range("S8").Value = 0.07793245400591471
$ws.Range("T8").Value = 0.07793245400591474

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 246.1811623333333
$ws.Range("H9").Value = 738.543487
$ws.Range("I9").Value = 0.5539389593320749
$ws.Range("J9").Value = 0.5539389593320749
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 93.562673
$ws.Range("N9").Value = 280.688019
$ws.Range("O9").Value = 0.1933186106880956
$ws.Range("P9").Value = 0.1933186106880956
$ws.Range("Q9").Value = 23033.36759015359
$ws.Range("R9").Value = 207300.3083113823
$ws.Range("S9").Value = 0.1070867100240862
$ws.Range("T9").Value = 0.1070867100240862

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.2401933333333333
$ws.Range("H10").Value = 0.72058
$ws.Range("I10").Value = 0.0005404655817044752
$ws.Range("J10").Value = 0.0005404655817044752
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 153.5290173333333
$ws.Range("N10").Value = 460.587052
$ws.Range("O10").Value = 0.3172206968818489
$ws.Range("P10").Value = 0.317220696881849
$ws.Range("Q10").Value = 36.87664643668444
$ws.Range("R10").Value = 331.88981793016
$ws.Range("S10").Value = 0.0001714468684689475
$ws.Range("T10").Value = 0.0001714468684689475

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.2401933333333333
$ws.Range("H11").Value = 0.72058
$ws.Range("I11").Value = 0.0005404655817044752
$ws.Range("J11").Value = 0.0005404655817044752
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 168.7997026666667
$ws.Range("N11").Value = 506.3991080000001
$ws.Range("O11").Value = 0.3487728915577651
$ws.Range("P11").Value = 0.3487728915577651
$ws.Range("Q11").Value = 40.54456324918223
$ws.Range("R11").Value = 364.9010692426401
$ws.Range("S11").Value = 0.0001884997437185194
$ws.Range("T11").Value = 0.0001884997437185194

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.2401933333333333
$ws.Range("H12").Value = 0.72058
$ws.Range("I12").Value = 0.0005404655817044752
$ws.Range("J12").Value = 0.0005404655817044752
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 68.09032333333333
$ws.Range("N12").Value = 204.27097
$ws.Range("O12").Value = 0.1406878008722904
$ws.Range("P12").Value = 0.1406878008722904
$ws.Range("Q12").Value = 16.35484172917778
$ws.Range("R12").Value = 147.1935755626
$ws.Range("S12").Value = 0.0000760369141371658
$ws.Range("T12").Value = 0.00007603691413716581

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.2401933333333333
$ws.Range("H13").Value = 0.72058
$ws.Range("I13").Value = 0.0005404655817044752
$ws.Range("J13").Value = 0.0005404655817044752
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 93.562673
$ws.Range("N13").Value = 280.688019
$ws.Range("O13").Value = 0.1933186106880956
$ws.Range("P13").Value = 0.1933186106880956
$ws.Range("Q13").Value = 22.47313030344667
$ws.Range("R13").Value = 202.25817273102
$ws.Range("S13").Value = 0.0001044820553798426
$ws.Range("T13").Value = 0.0001044820553798426

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 42.784114
$ws.Range("H14").Value = 128.352342
$ws.Range("I14").Value = 0.09626970382492123
$ws.Range("J14").Value = 0.09626970382492124
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 153.5290173333333
$ws.Range("N14").Value = 460.587052
$ws.Range("O14").Value = 0.3172206968818489
$ws.Range("P14").Value = 0.317220696881849
$ws.Range("Q14").Value = 6568.602979897308
$ws.Range("R14").Value = 59117.42681907578
$ws.Range("S14").Value = 0.0305387425359507
$ws.Range("T14").Value = 0.03053874253595072

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 42.784114
$ws.Range("H15").Value = 128.352342
$ws.Range("I15").Value = 0.09626970382492123
$ws.Range("J15").Value = 0.09626970382492124
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 168.7997026666667
$ws.Range("N15").Value = 506.3991080000001
$ws.Range("O15").Value = 0.3487728915577651
$ws.Range("P15").Value = 0.3487728915577651
$ws.Range("Q15").Value = 7221.945722056771
$ws.Range("R15").Value = 64997.51149851094
$ws.Range("S15").Value = 0.03357626297242741
$ws.Range("T15").Value = 0.03357626297242742

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 42.784114
$ws.Range("H16").Value = 128.352342
$ws.Range("I16").Value = 0.09626970382492123
$ws.Range("J16").Value = 0.09626970382492124
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 68.09032333333333
$ws.Range("N16").Value = 204.27097
$ws.Range("O16").Value = 0.1406878008722904
$ws.Range("P16").Value = 0.1406878008722904
$ws.Range("Q16").Value = 2913.184155790193
$ws.Range("R16").Value = 26218.65740211174
$ws.Range("S16").Value = 0.01354397292175489
$ws.Range("T16").Value = 0.01354397292175489

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 42.784114
$ws.Range("H17").Value = 128.352342
$ws.Range("I17").Value = 0.09626970382492123
$ws.Range("J17").Value = 0.09626970382492124
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 93.562673
$ws.Range("N17").Value = 280.688019
$ws.Range("O17").Value = 0.1933186106880956
$ws.Range("P17").Value = 0.1933186106880956
$ws.Range("Q17").Value = 4002.996067776722
$ws.Range("R17").Value = 36026.96460999049
$ws.Range("S17").Value = 0.01861072539478822
$ws.Range("T17").Value = 0.01861072539478822
